# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a different statistic (Strike#);
# this regenerates it with the correct per-game K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 3
    11 = 2
    12 = 1
    13 = 4
    14 = 3
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 2
    20 = 3
    21 = 2
    22 = 1
    23 = 6
    24 = 1
    25 = 2
    26 = 2
    27 = 1
    28 = 1
    29 = 2
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    36 = 1
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 0
    45 = 0
    47 = 1
    48 = 1
    49 = 3
    50 = 1
    51 = 2
    52 = 2
    53 = 1
    54 = 1
    55 = 2
    56 = 1
    57 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
